# Row 1, columns D:F used to hold real Excel date serial numbers
# (displayed through the "YYYY-MM-DD" number format already applied to
# those cells). The fix turns them into plain text labels showing the
# date (e.g. "2025.12.01") so the exported file name can be derived
# from these header values.
#
# Typing a string like "2025.12.01" straight into a date-formatted cell
# would make Excel re-parse it back into a date serial value, so the
# cells are temporarily switched to a "Text" number format before the
# text is assigned, and then the original "YYYY-MM-DD" number format
# (and therefore the same cell style) is restored. Only the stored
# value + type change (number -> shared string), exactly like the
# recorded edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateHeaderRange = $ws.Range("D1:F1")
$originalNumberFormat = $dateHeaderRange.NumberFormat

$dateHeaderRange.NumberFormat = "@"
$ws.Range("D1").Value = "2025.12.01"
$ws.Range("E1").Value = "2025.12.08"
$ws.Range("F1").Value = "2025.12.15"
$dateHeaderRange.NumberFormat = $originalNumberFormat
